$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain-numeric-looking text (e.g. "587.52") must be
# force-formatted as Text first, otherwise Excel auto-converts the assigned
# string into a floating-point number (losing the original text cell type).

$ws.Range("D2").Value = "69.420.38"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "3.401.93"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.52"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.35"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +11.53%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.49"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "689.34"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "3.955.18"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "69.525.78"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.426.69"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.121"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.38"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.35"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.52"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.07"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "558.15"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.65"
$ws.Range("E33").Value = "  +11.32%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "3.662.87"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.98"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.141"
$ws.Range("E39").Value = "  +5.06%  "
$ws.Range("D40").Value = "0.0₃0738"
$ws.Range("E40").Value = "  +9.46%  "
$ws.Range("E41").Value = "  +4.04%  "
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("E43").Value = "  +5.71%  "
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.09"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.72"
$ws.Range("E51").Value = "  +1.96%  "
